$wb = $excel.ActiveWorkbook

# Helper: write a value to a cell, preserving it as TEXT even when the
# string looks like a number (Excel would otherwise silently convert
# "-5.0" -> -5, "0.34" -> 0.34 (number), etc.). We briefly mark the cell
# as Text before assigning, then clear the formatting again so no visible
# number-format/style is left behind on the cell.
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$wsLider = $wb.Worksheets.Item("Restricciones_del_lider")
Set-TextValue $wsLider.Range("A2") "4.5 - x"
Set-TextValue $wsLider.Range("B2") "-5.0"
Set-TextValue $wsLider.Range("D2") "0.34"
Set-TextValue $wsLider.Range("A3") "-4.5 + x"
Set-TextValue $wsLider.Range("B3") "4.0"
Set-TextValue $wsLider.Range("D3") "0.0"

$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")
Set-TextValue $wsFollower.Range("A2") "2.8 - y"
Set-TextValue $wsFollower.Range("B2") "-3.8"
Set-TextValue $wsFollower.Range("D2") "0.09"
Set-TextValue $wsFollower.Range("E2") "5.6000000000000005"
Set-TextValue $wsFollower.Range("F2") "5.6000000000000005"
Set-TextValue $wsFollower.Range("A3") "-2.8 + y"
Set-TextValue $wsFollower.Range("B3") "1.7999999999999998"
Set-TextValue $wsFollower.Range("D3") "0.82"
Set-TextValue $wsFollower.Range("E3") "8.299999999999999"
Set-TextValue $wsFollower.Range("F3") "7.199999999999999"

$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto.Range("A2") "4.5"
Set-TextValue $wsPunto.Range("B2") "2.8"

$wsVecbf = $wb.Worksheets.Item("Vector_bf")
Set-TextValue $wsVecbf.Range("A2") "-5.23"

$wsVecBF = $wb.Worksheets.Item("Vector_BF")
Set-TextValue $wsVecBF.Range("A2") "1.34"
Set-TextValue $wsVecBF.Range("A3") "-1.6999999999999984"
